# Refactor a tarea de cargue
# The "opcionSubmenu" column (values "Recargar") in the Datos sheet is no
# longer needed by the loading task, so remove it entirely and let the
# following columns (tipoCuenta / numeroCuenta) shift left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Column N holds the header "opcionSubmenu" and the value "Recargar" for
# every data row - delete the whole column so O/P (tipoCuenta/numeroCuenta)
# shift into N/O.
$ws.Range("N1").EntireColumn.Delete()

# Reflect the resulting selection/zoom state of the sheet.
[void]$ws.Range("N12").Select()
$excel.ActiveWindow.Zoom = 164
